$wb = $excel.ActiveWorkbook

# Update the "Ready for handoff" -> "In Translation" status on each sheet.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Adjust column widths (E,F on Overview; C on zh-cn and de-de)
$overview.Columns.Item(5).ColumnWidth = 13.4101845877511
$overview.Columns.Item(6).ColumnWidth = 13.4101845877511

$zhcn.Columns.Item(3).ColumnWidth = 13.4101845877511
$dede.Columns.Item(3).ColumnWidth = 13.4101845877511
